$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price values are pure numeric-looking strings (e.g. '492.33').
# Excel's COM Range.Value setter would auto-convert these to numbers (losing the
# original text formatting / exact decimal text), so force those specific cells to
# Text format first, matching how the sheet's author kept them as text.
$textCells = @(
    "D5", "D6", "D8", "D10", "D11", "D12", "D16", "D19",
    "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27",
    "D28", "D29", "D33", "D34", "D35", "D36", "D37", "D39",
    "D40", "D42", "D43", "D44", "D45", "D47", "D48", "D49",
    "D50", "D51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated coin data (rank order, swapped rows, prices, and 1h volume %).
$ws.Range("D2").Value = '57.689.62'
$ws.Range("E2").Value = '  -4.98%  '
$ws.Range("D3").Value = '2.686.16'
$ws.Range("E3").Value = '  -7.72%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '492.33'
$ws.Range("E5").Value = '  -6.89%  '
$ws.Range("D6").Value = '133.72'
$ws.Range("E6").Value = '  -6.82%  '
$ws.Range("E7").Value = '  +0.45%  '
$ws.Range("D8").Value = '0.518'
$ws.Range("E8").Value = '  -6.67%  '
$ws.Range("D9").Value = '2.651.08'
$ws.Range("E9").Value = '  -9.01%  '
$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").Value = '0.100'
$ws.Range("E10").Value = '  -7.24%  '
$ws.Range("B11").Value = 'Toncoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D11").Value = '5.79'
$ws.Range("E11").Value = '  -1.42%  '
$ws.Range("D12").Value = '0.337'
$ws.Range("E12").Value = '  -4.48%  '
$ws.Range("E13").Value = '  +0.94%  '
$ws.Range("D14").Value = '3.176.68'
$ws.Range("E14").Value = '  -7.03%  '
$ws.Range("D15").Value = '58.010.95'
$ws.Range("E15").Value = '  -4.62%  '
$ws.Range("D16").Value = '20.86'
$ws.Range("E16").Value = '  -8.09%  '
$ws.Range("D17").Value = '2.702.78'
$ws.Range("E17").Value = '  -7.25%  '
$ws.Range("E18").Value = '  -6.60%  '
$ws.Range("D19").Value = '4.58'
$ws.Range("E19").Value = '  -7.03%  '
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").Value = '337.99'
$ws.Range("E20").Value = '  -6.45%  '
$ws.Range("B21").Value = 'Chainlink'
$ws.Range("C21").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D21").Value = '10.60'
$ws.Range("E21").Value = '  -8.33%  '
$ws.Range("D22").Value = '6.06'
$ws.Range("E22").Value = '  -7.54%  '
$ws.Range("D23").Value = '0.996'
$ws.Range("E23").Value = '  -0.45%  '
$ws.Range("D24").Value = '5.60'
$ws.Range("E24").Value = '  -1.37%  '
$ws.Range("D25").Value = '61.29'
$ws.Range("E25").Value = '  -3.46%  '
$ws.Range("D26").Value = '0.413'
$ws.Range("E26").Value = '  -8.51%  '
$ws.Range("D27").Value = '0.167'
$ws.Range("E27").Value = '  -8.41%  '
$ws.Range("D28").Value = '1.01'
$ws.Range("E28").Value = '  +0.73%  '
$ws.Range("D29").Value = '7.18'
$ws.Range("E29").Value = '  -6.75%  '
$ws.Range("D30").Value = '0.0₃0786'
$ws.Range("E30").Value = '  -8.87%  '
$ws.Range("E31").Value = '  +0.14%  '
$ws.Range("E32").Value = '  -6.35%  '
$ws.Range("D33").Value = '18.51'
$ws.Range("E33").Value = '  -6.19%  '
$ws.Range("D34").Value = '145.54'
$ws.Range("E34").Value = '  -6.04%  '
$ws.Range("D35").Value = '4.05'
$ws.Range("E35").Value = '  -6.97%  '
$ws.Range("D36").Value = '5.16'
$ws.Range("E36").Value = '  -7.77%  '
$ws.Range("D37").Value = '0.892'
$ws.Range("E37").Value = '  -11.31%  '
$ws.Range("E38").Value = '  -9.46%  '
$ws.Range("D39").Value = '35.68'
$ws.Range("E39").Value = '  -6.01%  '
$ws.Range("D40").Value = '0.996'
$ws.Range("E40").Value = '  -0.32%  '
$ws.Range("D41").Value = '2.126.99'
$ws.Range("E41").Value = '  -8.83%  '
$ws.Range("D42").Value = '3.42'
$ws.Range("E42").Value = '  -7.40%  '
$ws.Range("B43").Value = 'Stacks'
$ws.Range("C43").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D43").Value = '1.32'
$ws.Range("E43").Value = '  -10.84%  '
$ws.Range("B44").Value = 'Hedera'
$ws.Range("C44").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D44").Value = '0.0541'
$ws.Range("E44").Value = '  -5.22%  '
$ws.Range("D45").Value = '0.575'
$ws.Range("E45").Value = '  -10.67%  '
$ws.Range("E46").Value = '  +0.12%  '
$ws.Range("D47").Value = '18.44'
$ws.Range("E47").Value = '  -11.81%  '
$ws.Range("D48").Value = '0.0221'
$ws.Range("E48").Value = '  -5.24%  '
$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D49").Value = '4.49'
$ws.Range("E49").Value = '  -8.41%  '
$ws.Range("B50").Value = 'Stellar'
$ws.Range("C50").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D50").Value = '0.0871'
$ws.Range("E50").Value = '  -5.93%  '
$ws.Range("D51").Value = '17.07'
$ws.Range("E51").Value = '  -7.53%  '
